# Adds a description-of-attributes sub-table below the existing ER table
# on sheet1: a merged/centered "Cliente" header row, then an
# Atributo/Tipo/Obligatorio/llave prim header row, followed by two data
# rows describing the "ID" and "cuentaPlatzi" attributes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers for the attribute table (written first so the new shared
# strings land in the same order the original authored workbook used)
$ws.Range("A10").Value = "Atributo"
$ws.Range("B10").Value = "Tipo"
$ws.Range("C10").Value = "Obligatorio"
$ws.Range("D10").Value = "llave prim"

# Header: merged & centered "Cliente" across A9:D9
$ws.Range("A9").Value = "Cliente"
$ws.Range("A9").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A9:D9").Merge() | Out-Null

# Row describing the ID attribute (primary key, required)
$ws.Range("A11").Value = "ID"
$ws.Range("B11").Value = "Serial"
$ws.Range("C11").Value = "x"
$ws.Range("D11").Value = "x"

# Row describing the cuentaPlatzi attribute
$ws.Range("A12").Value = "cuentaPlatzi"
$ws.Range("B12").Value = "varchar(60)"
$ws.Range("C12").Value = "x"

# Reposition the view/selection similar to the authored edit (scrolls the
# viewport so row 6 is at the top, then selects D12 as the active cell)
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("D12").Select()
